# Updates crypto list values (price/volume/name/link) per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '38.661.40'
$ws.Range("E2").Value = '  -5.09%  '

# Row 3
$ws.Range("D3").Value = '2.178.70'

# Row 4
$ws.Range("E4").Value = '  +0.46%  '

# Row 5
$ws.Range("D5").Value = '''293.92'
$ws.Range("E5").Value = '  -5.65%  '

# Row 6
$ws.Range("D6").Value = '''79.91'
$ws.Range("E6").Value = '  -8.92%  '

# Row 7
$ws.Range("D7").Value = '''0.500'
$ws.Range("E7").Value = '  -5.39%  '

# Row 8
$ws.Range("E8").Value = '  +0.20%  '

# Row 9
$ws.Range("D9").Value = '''0.454'
$ws.Range("E9").Value = '  -8.47%  '

# Row 10
$ws.Range("D10").Value = '''0.0766'
$ws.Range("E10").Value = '  -8.75%  '

# Row 11
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").Value = '''27.70'
$ws.Range("E11").Value = '  -10.43%  '

# Row 12
$ws.Range("B12").Value = 'OKB'
$ws.Range("C12").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D12").Value = '''45.63'
$ws.Range("E12").Value = '  -13.35%  '

# Row 13
$ws.Range("D13").Value = '''0.107'
$ws.Range("E13").Value = '  -2.55%  '

# Row 14
$ws.Range("D14").Value = '2.530.14'
$ws.Range("E14").Value = '  -7.58%  '

# Row 15
$ws.Range("D15").Value = '''6.06'
$ws.Range("E15").Value = '  -7.99%  '

# Row 16
$ws.Range("D16").Value = '''13.72'
$ws.Range("E16").Value = '  -8.97%  '

# Row 17
$ws.Range("D17").Value = '2.197.91'
$ws.Range("E17").Value = '  -6.90%  '

# Row 18
$ws.Range("D18").Value = '''0.699'
$ws.Range("E18").Value = '  -8.68%  '

# Row 19
$ws.Range("D19").Value = '38.632.03'
$ws.Range("E19").Value = '  -4.81%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0855'
$ws.Range("E20").Value = '  -6.48%  '

# Row 21
$ws.Range("D21").Value = '''5.60'
$ws.Range("E21").Value = '  -9.15%  '

# Row 22
$ws.Range("D22").Value = '''63.59'
$ws.Range("E22").Value = '  -7.69%  '

# Row 23
$ws.Range("D23").Value = '''9.72'
$ws.Range("E23").Value = '  -11.10%  '

# Row 24
$ws.Range("D24").Value = '''221.16'
$ws.Range("E24").Value = '  -5.35%  '

# Row 25
$ws.Range("E25").Value = '  +0.10%  '

# Row 26
$ws.Range("D26").Value = '''2.35'
$ws.Range("E26").Value = '  -10.94%  '

# Row 27
$ws.Range("D27").Value = '''1.72'
$ws.Range("E27").Value = '  -5.83%  '

# Row 28
$ws.Range("D28").Value = '''22.06'
$ws.Range("E28").Value = '  -7.78%  '

# Row 29
$ws.Range("E29").Value = '  -2.40%  '

# Row 30
$ws.Range("D30").Value = '''8.83'
$ws.Range("E30").Value = '  -6.23%  '

# Row 31
$ws.Range("D31").Value = '''146.80'
$ws.Range("E31").Value = '  -4.04%  '

# Row 32
$ws.Range("D32").Value = '''30.79'
$ws.Range("E32").Value = '  -9.43%  '

# Row 33
$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  -0.05%  '

# Row 34
$ws.Range("D34").Value = '''4.69'
$ws.Range("E34").Value = '  -10.87%  '

# Row 35
$ws.Range("E35").Value = '  -5.87%  '

# Row 36
$ws.Range("D36").Value = '''0.0680'
$ws.Range("E36").Value = '  -7.33%  '

# Row 37
$ws.Range("D37").Value = '''0.108'
$ws.Range("E37").Value = '  -5.13%  '

# Row 38
$ws.Range("D38").Value = '''0.0940'
$ws.Range("E38").Value = '  -5.90%  '

# Row 39
$ws.Range("D39").Value = '''2.58'
$ws.Range("E39").Value = '  -7.56%  '

# Row 40
$ws.Range("D40").Value = '''1.57'
$ws.Range("E40").Value = '  -9.08%  '

# Row 41
$ws.Range("D41").Value = '''14.06'
$ws.Range("E41").Value = '  -12.83%  '

# Row 42
$ws.Range("D42").Value = '''3.58'
$ws.Range("E42").Value = '  -7.79%  '

# Row 43
$ws.Range("D43").Value = '1.872.00'
$ws.Range("E43").Value = '  -4.66%  '

# Row 44
$ws.Range("D44").Value = '''2.02'
$ws.Range("E44").Value = '  -15.45%  '

# Row 45
$ws.Range("D45").Value = '''0.0252'
$ws.Range("E45").Value = '  -6.92%  '

# Row 46
$ws.Range("D46").Value = '''15.82'
$ws.Range("E46").Value = '  -10.65%  '

# Row 47
$ws.Range("D47").Value = '''8.81'
$ws.Range("E47").Value = '  -8.77%  '

# Row 48
$ws.Range("D48").Value = '''2.51'
$ws.Range("E48").Value = '  -8.44%  '

# Row 49
$ws.Range("D49").Value = '2.383.64'

# Row 50
$ws.Range("D50").Value = '''69.18'
$ws.Range("E50").Value = '  -4.81%  '

# Row 51
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = '''1.05'
$ws.Range("E51").Value = '  -2.07%  '
